$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.039021124875183
$ws.Range("D2").Value2 = 1.047652440371431
$ws.Range("E2").Value2 = 1.048044455949791
$ws.Range("F2").Value2 = 1.059970098613571
$ws.Range("I2").Value2 = 1.044724882390172
$ws.Range("J2").Value2 = 1.044115400443574
$ws.Range("K2").Value2 = 1.050414479020055
$ws.Range("L2").Value2 = 1.050805399678235
$ws.Range("M2").Value2 = 1.062698152631818
$ws.Range("N2").Value2 = 1.018640176534435
$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.039786323655679
$ws.Range("D3").Value2 = 1.048245561182424
$ws.Range("E3").Value2 = 1.048713355735134
$ws.Range("F3").Value2 = 1.060692375683804
$ws.Range("I3").Value2 = 1.044914257711019
$ws.Range("J3").Value2 = 1.044526556064454
$ws.Range("K3").Value2 = 1.050820054402176
$ws.Range("L3").Value2 = 1.051286636303312
$ws.Range("M3").Value2 = 1.063234996180676
$ws.Range("N3").Value2 = 1.018777007648995
$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.040282170622021
$ws.Range("D4").Value2 = 1.04862994224295
$ws.Range("E4").Value2 = 1.049147175216379
$ws.Range("F4").Value2 = 1.061160743243147
$ws.Range("I4").Value2 = 1.045035848890331
$ws.Range("J4").Value2 = 1.044792610248827
$ws.Range("K4").Value2 = 1.051082381177165
$ws.Range("L4").Value2 = 1.051598340488279
$ws.Range("M4").Value2 = 1.063582701084819
$ws.Range("N4").Value2 = 1.018865530562861
$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.040490793192375
$ws.Range("D5").Value2 = 1.048791675825067
$ws.Range("E5").Value2 = 1.049329789425624
$ws.Range("F5").Value2 = 1.061357883575054
$ws.Range("I5").Value2 = 1.045086738105743
$ws.Range("J5").Value2 = 1.044904460194156
$ws.Range("K5").Value2 = 1.051192636131186
$ws.Range("L5").Value2 = 1.051729454007581
$ws.Range("M5").Value2 = 1.06372895373729
$ws.Range("N5").Value2 = 1.0189027412812
$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.040525831663231
$ws.Range("D6").Value2 = 1.048818839744013
$ws.Range("E6").Value2 = 1.049360464962581
$ws.Range("F6").Value2 = 1.061390998231597
$ws.Range("I6").Value2 = 1.045095269248604
$ws.Range("J6").Value2 = 1.044923240296648
$ws.Range("K6").Value2 = 1.051211146784356
$ws.Range("L6").Value2 = 1.051751472783519
$ws.Range("M6").Value2 = 1.063753514694464
$ws.Range("N6").Value2 = 1.018908988857208
$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.040284957587403
$ws.Range("D7").Value2 = 1.048632102788251
$ws.Range("E7").Value2 = 1.049149614389729
$ws.Range("F7").Value2 = 1.061163376507143
$ws.Range("I7").Value2 = 1.045036529770445
$ws.Range("J7").Value2 = 1.044794104792006
$ws.Range("K7").Value2 = 1.051083854518054
$ws.Range("L7").Value2 = 1.05160009214839
$ws.Range("M7").Value2 = 1.063584655017602
$ws.Range("N7").Value2 = 1.018866027791874
$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.039279578723622
$ws.Range("D8").Value2 = 1.047852764500149
$ws.Range("E8").Value2 = 1.048270306159387
$ws.Range("F8").Value2 = 1.060213985894119
$ws.Range("I8").Value2 = 1.044789078257121
$ws.Range("J8").Value2 = 1.044254349711177
$ws.Range("K8").Value2 = 1.050551566525803
$ws.Range("L8").Value2 = 1.050967970309539
$ws.Range("M8").Value2 = 1.062879511765779
$ws.Range("N8").Value2 = 1.018686422260593
$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.037513515466408
$ws.Range("D9").Value2 = 1.046484094962631
$ws.Range("E9").Value2 = 1.046728577651718
$ws.Range("F9").Value2 = 1.058548841563838
$ws.Range("I9").Value2 = 1.044345822800515
$ws.Range("J9").Value2 = 1.043303362700359
$ws.Range("K9").Value2 = 1.049612847335519
$ws.Range("L9").Value2 = 1.049856550084739
$ws.Range("M9").Value2 = 1.061639575070342
$ws.Range("N9").Value2 = 1.018369833831367
$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.036339985884746
$ws.Range("D10").Value2 = 1.045574877063405
$ws.Range("E10").Value2 = 1.045706077306734
$ws.Range("F10").Value2 = 1.057444125101401
$ws.Range("I10").Value2 = 1.044045524194373
$ws.Range("J10").Value2 = 1.042669542133062
$ws.Range("K10").Value2 = 1.048986610959007
$ws.Range("L10").Value2 = 1.049117351915206
$ws.Range("M10").Value2 = 1.060814816195663
$ws.Range("N10").Value2 = 1.018158736701915
$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.03583277208531
$ws.Range("D11").Value2 = 1.045181966778111
$ws.Range("E11").Value2 = 1.045264611391704
$ws.Range("F11").Value2 = 1.056967074983664
$ws.Range("I11").Value2 = 1.043914366027329
$ws.Range("J11").Value2 = 1.042395148946265
$ws.Range("K11").Value2 = 1.048715361715556
$ws.Range("L11").Value2 = 1.048797705897433
$ws.Range("M11").Value2 = 1.060458151660875
$ws.Range("N11").Value2 = 1.018067326494092
$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.035644512407795
$ws.Range("D12").Value2 = 1.045036142719993
$ws.Range("E12").Value2 = 1.045100826290834
$ws.Range("F12").Value2 = 1.056790074861964
$ws.Range("I12").Value2 = 1.043865479680682
$ws.Range("J12").Value2 = 1.042293236835717
$ws.Range("K12").Value2 = 1.048614596433654
$ws.Range("L12").Value2 = 1.048679041630871
$ws.Range("M12").Value2 = 1.0603257417524
$ws.Range("N12").Value2 = 1.018033372613955
$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.035684888283916
$ws.Range("D13").Value2 = 1.045067416997745
$ws.Range("E13").Value2 = 1.045135949892203
$ws.Range("F13").Value2 = 1.056828033013053
$ws.Range("I13").Value2 = 1.043875973571713
$ws.Range("J13").Value2 = 1.042315096865874
$ws.Range("K13").Value2 = 1.048636211413064
$ws.Range("L13").Value2 = 1.048704492483388
$ws.Range("M13").Value2 = 1.060354140876594
$ws.Range("N13").Value2 = 1.018040655831059
$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.035817207570749
$ws.Range("D14").Value2 = 1.045169910449617
$ws.Range("E14").Value2 = 1.04525106887379
$ws.Range("F14").Value2 = 1.056952440055424
$ws.Range("I14").Value2 = 1.043910328495799
$ws.Range("J14").Value2 = 1.042386724652544
$ws.Range("K14").Value2 = 1.04870703264646
$ws.Range("L14").Value2 = 1.048787895706134
$ws.Range("M14").Value2 = 1.06044720515382
$ws.Range("N14").Value2 = 1.018064519854048
$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.035898752742931
$ws.Range("D15").Value2 = 1.04523307605288
$ws.Range("E15").Value2 = 1.045322023381402
$ws.Range("F15").Value2 = 1.057029117584272
$ws.Range("I15").Value2 = 1.043931473416668
$ws.Range("J15").Value2 = 1.042430858228777
$ws.Range("K15").Value2 = 1.048750666501733
$ws.Range("L15").Value2 = 1.048839292044891
$ws.Range("M15").Value2 = 1.060504554613737
$ws.Range("N15").Value2 = 1.018079223273527
$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.036373667810691
$ws.Range("D16").Value2 = 1.045600969984252
$ws.Range("E16").Value2 = 1.045735403172197
$ws.Range("F16").Value2 = 1.057475812920947
$ws.Range("I16").Value2 = 1.044054205067817
$ws.Range("J16").Value2 = 1.04268775399224
$ws.Range("K16").Value2 = 1.049004611220135
$ws.Range("L16").Value2 = 1.049138575003764
$ws.Range("M16").Value2 = 1.060838496727028
$ws.Range("N16").Value2 = 1.018164803260365
$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.036671820523793
$ws.Range("D17").Value2 = 1.045831952309083
$ws.Range("E17").Value2 = 1.045995050549694
$ws.Range("F17").Value2 = 1.057756362624217
$ws.Range("I17").Value2 = 1.044130890352298
$ws.Range("J17").Value2 = 1.042848913731827
$ws.Range("K17").Value2 = 1.049163882398769
$ws.Range("L17").Value2 = 1.049326424027571
$ws.Range("M17").Value2 = 1.061048094652149
$ws.Range("N17").Value2 = 1.018218484676265
$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.036845817768385
$ws.Range("D18").Value2 = 1.045966756193124
$ws.Range("E18").Value2 = 1.046146622164358
$ws.Range("F18").Value2 = 1.057920127738055
$ws.Range("I18").Value2 = 1.044175510789883
$ws.Range("J18").Value2 = 1.042942920721423
$ws.Range("K18").Value2 = 1.0492567743041
$ws.Range("L18").Value2 = 1.049436034724586
$ws.Range("M18").Value2 = 1.061170393980281
$ws.Range("N18").Value2 = 1.018249795741509
$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.036905161503488
$ws.Range("D19").Value2 = 1.046012733628255
$ws.Range("E19").Value2 = 1.046198325082339
$ws.Range("F19").Value2 = 1.057975988578614
$ws.Range("I19").Value2 = 1.044190706724542
$ws.Range("J19").Value2 = 1.042974975520864
$ws.Range("K19").Value2 = 1.049288446593476
$ws.Range("L19").Value2 = 1.04947341613547
$ws.Range("M19").Value2 = 1.061212102366283
$ws.Range("N19").Value2 = 1.018260471922969
$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.036639822254971
$ws.Range("D20").Value2 = 1.045807162236163
$ws.Range("E20").Value2 = 1.045967180033668
$ws.Range("F20").Value2 = 1.057726249327548
$ws.Range("I20").Value2 = 1.04412267399306
$ws.Range("J20").Value2 = 1.042831622264667
$ws.Range("K20").Value2 = 1.049146794954108
$ws.Range("L20").Value2 = 1.049306265296619
$ws.Range("M20").Value2 = 1.061025602172128
$ws.Range("N20").Value2 = 1.018212725205734
$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.035778238899309
$ws.Range("D21").Value2 = 1.04513972534472
$ws.Range("E21").Value2 = 1.045217163802742
$ws.Range("F21").Value2 = 1.056915799810791
$ws.Range("I21").Value2 = 1.043900216467047
$ws.Range("J21").Value2 = 1.042365631774414
$ws.Range("K21").Value2 = 1.048686177856527
$ws.Range("L21").Value2 = 1.048763333686987
$ws.Range("M21").Value2 = 1.060419798064018
$ws.Range("N21").Value2 = 1.018057492491643
$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.035237350475171
$ws.Range("D22").Value2 = 1.044720778413431
$ws.Range("E22").Value2 = 1.044746728088125
$ws.Range("F22").Value2 = 1.056407382037782
$ws.Range("I22").Value2 = 1.043759374956589
$ws.Range("J22").Value2 = 1.042072702095548
$ws.Range("K22").Value2 = 1.048396505340748
$ws.Range("L22").Value2 = 1.048422356576052
$ws.Range("M22").Value2 = 1.060039317715065
$ws.Range("N22").Value2 = 1.017959891431966
$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.03552400694532
$ws.Range("D23").Value2 = 1.044942803338784
$ws.Range("E23").Value2 = 1.044996007198162
$ws.Range("F23").Value2 = 1.056676794709322
$ws.Range("I23").Value2 = 1.04383412968172
$ws.Range("J23").Value2 = 1.042227983777805
$ws.Range("K23").Value2 = 1.04855007181374
$ws.Range("L23").Value2 = 1.048603077901115
$ws.Range("M23").Value2 = 1.060240977808827
$ws.Range("N23").Value2 = 1.018011631440389
$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.036654280618073
$ws.Range("D24").Value2 = 1.045818363569106
$ws.Range("E24").Value2 = 1.045979773138297
$ws.Range("F24").Value2 = 1.057739855843513
$ws.Range("I24").Value2 = 1.044126386948401
$ws.Range("J24").Value2 = 1.042839435518374
$ws.Range("K24").Value2 = 1.04915451606061
$ws.Range("L24").Value2 = 1.049315374031286
$ws.Range("M24").Value2 = 1.061035765418719
$ws.Range("N24").Value2 = 1.018215327663846
$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.037969416017478
$ws.Range("D25").Value2 = 1.046837368539433
$ws.Range("E25").Value2 = 1.04712622310206
$ws.Range("F25").Value2 = 1.05897838281434
$ws.Range("I25").Value2 = 1.044461264173919
$ws.Range("J25").Value2 = 1.043549191617121
$ws.Range("K25").Value2 = 1.04985560978693
$ws.Range("L25").Value2 = 1.050143577604928
$ws.Range("M25").Value2 = 1.061959807419966
$ws.Range("N25").Value2 = 1.018451688408978
